$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

$d = $word.ActiveDocument

$d.Content.Find.Execute("p for child in children", $true, $false, $false, $false, $false, $true, 1, $false, "p for child in children.filter(legal_parent=`"both`")", 2)
